$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cells value while forcing it to stay literal Text (the "Price"
# column holds strings like "320.59" / "1.00" that Excel would otherwise silently
# coerce to numbers, dropping the original formatting). We flip the cell to the
# Text number format just long enough to set the value, then restore the default
# "Normal" style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "47.265.62"
$ws.Range("E2").Value = "  -0.31%  "

Set-TextValue "D3" "2.488.78"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue "D5" "320.59"
$ws.Range("E5").Value = "  -1.24%  "

Set-TextValue "D6" "107.87"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("E7").Value = "  -0.68%  "

Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.02%  "

Set-TextValue "D9" "0.534"
$ws.Range("E9").Value = "  -1.79%  "

Set-TextValue "D10" "38.60"
$ws.Range("E10").Value = "  +5.17%  "

Set-TextValue "D11" "0.0808"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("E12").Value = "  -0.02%  "

Set-TextValue "D13" "18.20"
$ws.Range("E13").Value = "  -0.86%  "

Set-TextValue "D14" "7.12"
$ws.Range("E14").Value = "  -0.85%  "

Set-TextValue "D15" "2.856.62"
$ws.Range("E15").Value = "  -1.67%  "

Set-TextValue "D16" "2.472.25"
$ws.Range("E16").Value = "  -2.58%  "

Set-TextValue "D17" "0.846"
$ws.Range("E17").Value = "  -0.36%  "

Set-TextValue "D18" "47.176.36"
$ws.Range("E18").Value = "  -0.28%  "

Set-TextValue "D19" "12.73"
$ws.Range("E19").Value = "  +0.50%  "

Set-TextValue "D20" "6.61"
$ws.Range("E20").Value = "  +1.53%  "

Set-TextValue "D21" "0.0₃0932"
$ws.Range("E21").Value = "  -1.08%  "

Set-TextValue "D22" "2.71"
$ws.Range("E22").Value = "  +13.24%  "

Set-TextValue "D23" "70.35"
$ws.Range("E23").Value = "  -0.93%  "

Set-TextValue "D24" "245.23"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("E26").Value = "  +0.05%  "

Set-TextValue "D27" "25.69"
$ws.Range("E27").Value = "  -2.98%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "2.28"
$ws.Range("E28").Value = "  +3.24%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "10.00"
$ws.Range("E29").Value = "  +1.08%  "

Set-TextValue "D30" "34.77"
$ws.Range("E30").Value = "  -1.48%  "

Set-TextValue "D31" "0.136"
$ws.Range("E31").Value = "  +0.16%  "

Set-TextValue "D32" "49.49"
$ws.Range("E32").Value = "  -0.48%  "

Set-TextValue "D33" "20.06"
$ws.Range("E33").Value = "  +1.55%  "

Set-TextValue "D34" "5.35"
$ws.Range("E34").Value = "  +0.33%  "

Set-TextValue "D35" "0.0780"
$ws.Range("E35").Value = "  +0.51%  "

$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +1.12%  "

$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "22.06"
$ws.Range("E42").Value = "  +5.66%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "119.38"
$ws.Range("E43").Value = "  -3.11%  "

Set-TextValue "D44" "0.0295"
$ws.Range("E44").Value = "  -0.53%  "

Set-TextValue "D45" "1.984.58"
$ws.Range("E45").Value = "  -0.08%  "

Set-TextValue "D46" "3.01"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("E47").Value = "  -5.97%  "

$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").Value = "  -3.44%  "

Set-TextValue "D50" "5.11"
$ws.Range("E50").Value = "  -6.54%  "

Set-TextValue "D51" "56.82"
$ws.Range("E51").Value = "  +3.55%  "
